# Adds measurement units to various simulation-parameter values in the
# MD protocol table (image/document extraction correction pass).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# dt values -> add "fs" (femtoseconds)
$ws.Range("C15").Value = "2 fs"
$ws.Range("C18").Value = "4 fs"
$ws.Range("C25").Value = "4 fs"

# cut values -> add "Å" (Angstrom)
$ws.Range("C19").Value = "8 Å"
$ws.Range("C26").Value = "8 Å"

# restraint_wt (Minimization) -> reformat exponent with carets
$ws.Range("C30").Value = "25 kcal mol-1 Å^-2^"

# Thermalization simulation time -> add "ps" (picoseconds)
$ws.Range("C32").Value = "50 ps"

# Equilibration (paragraph 13) values
$ws.Range("C35").Value = "300 ps"
$ws.Range("C36").Value = "1 atm"
$ws.Range("C37").Value = "300 K"
$ws.Range("C38").Value = "10 kcal mol^-1^ Å^-2^"

# Equilibration (paragraph 14) simulation time -> add "ns" (nanoseconds)
$ws.Range("C39").Value = "300 ns"

# Equilibration (paragraph 15) values
$ws.Range("C41").Value = "0 kcal mol^-1^ Å^-2^"
$ws.Range("C42").Value = "250 ns"
